$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the bold/bordered/centered format used in column A (e.g. A2) to the two new rows
# that extend the table (A30, A31) before we populate them.
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(30, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(31, 1).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Update A (index) and B (label) columns for rows 3-31.
# Two new rows were added near the top of the simulated data (new labels "Holden" and
# "Rizzie Spiral"), which pushed every later label down by two rows; "Thomas Hex" was also
# renamed to "Matthies Hex".
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = "Spiral5"
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = "Holden"
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = "Rizzie Spiral"
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = "RotRing OmegaMax-90"
$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = "Equal Angle"
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "Tilt Rotate"
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "CLR"
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "Rizzie Hex"
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "Matthies Hex"
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = "Tilt Rotate_Partial"
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "RotRing OmegaMax-60"
$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = "Equal Angle_Partial"
$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = "Rizzie Hex_Partial"
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "ND Single"
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "RD Single"
$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = "TD Single"
$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = "Morris Single"
$ws.Cells.Item(20, 1).Value = 18
$ws.Cells.Item(20, 2).Value = "Ring Perpendicular to ND"
$ws.Cells.Item(21, 1).Value = 19
$ws.Cells.Item(21, 2).Value = "Ring Perpendicular to RD"
$ws.Cells.Item(22, 1).Value = 20
$ws.Cells.Item(22, 2).Value = "Ring Perpendicular to TD"
$ws.Cells.Item(23, 1).Value = 21
$ws.Cells.Item(23, 2).Value = "OffsetFTD"
$ws.Cells.Item(24, 1).Value = 22
$ws.Cells.Item(24, 2).Value = "OffsetATD"
$ws.Cells.Item(25, 1).Value = 23
$ws.Cells.Item(25, 2).Value = "OffsetF45"
$ws.Cells.Item(26, 1).Value = 24
$ws.Cells.Item(26, 2).Value = "OffsetA45"
$ws.Cells.Item(27, 1).Value = 25
$ws.Cells.Item(27, 2).Value = "OffsetFRD"
$ws.Cells.Item(28, 1).Value = 26
$ws.Cells.Item(28, 2).Value = "OffsetARD"
$ws.Cells.Item(29, 1).Value = 27
$ws.Cells.Item(29, 2).Value = "Gaussian Quadrature"
$ws.Cells.Item(30, 1).Value = 28
$ws.Cells.Item(30, 2).Value = "Michael-CCHex"
$ws.Cells.Item(31, 1).Value = 29
$ws.Cells.Item(31, 2).Value = "Michael-SNHex"

# Update the simulated C:W data for rows 3-31 (re-ran simulation; two extra rows of
# results were produced at the top, pushing prior rows two down and appending two more
# at the bottom).
$ws.Cells.Item(3, 3).Value = 0.9999866499529582
$ws.Cells.Item(3, 4).Value = 0.9999332499417307
$ws.Cells.Item(3, 5).Value = 1.000044500320241
$ws.Cells.Item(3, 6).Value = 0.9999332499417307
$ws.Cells.Item(3, 7).Value = 1.000016686528213
$ws.Cells.Item(3, 8).Value = 1.000024777586765
$ws.Cells.Item(3, 9).Value = 0.9999856573744318
$ws.Cells.Item(3, 10).Value = 1.000044500320241
$ws.Cells.Item(3, 11).Value = 1.000044500320241
$ws.Cells.Item(3, 12).Value = 0.9999566000146138
$ws.Cells.Item(3, 13).Value = 1.000016686528213
$ws.Cells.Item(3, 14).Value = 1.000044500320241
$ws.Cells.Item(3, 15).Value = 1.000016686528213
$ws.Cells.Item(3, 16).Value = 0.9999749682349721
$ws.Cells.Item(3, 17).Value = 1.000001171951323
$ws.Cells.Item(3, 18).Value = 0.9999981455967283
$ws.Cells.Item(3, 19).Value = 0.9999785312814587
$ws.Cells.Item(3, 20).Value = 0.9999981455967283
$ws.Cells.Item(3, 21).Value = 0.9999950235411542
$ws.Cells.Item(3, 22).Value = 1.000004918896972
$ws.Cells.Item(3, 23).Value = 0.9999956010308959
$ws.Cells.Item(4, 3).Value = 0.9995277030292372
$ws.Cells.Item(4, 4).Value = 0.9976385145666206
$ws.Cells.Item(4, 5).Value = 1.001574323867359
$ws.Cells.Item(4, 6).Value = 0.9976385145666206
$ws.Cells.Item(4, 7).Value = 1.000590371197722
$ws.Cells.Item(4, 8).Value = 1.00087656457005
$ws.Cells.Item(4, 9).Value = 0.9994925735891693
$ws.Cells.Item(4, 10).Value = 1.001574323867359
$ws.Cells.Item(4, 11).Value = 1.001574323867359
$ws.Cells.Item(4, 12).Value = 0.9984645493474288
$ws.Cells.Item(4, 13).Value = 1.000590371197722
$ws.Cells.Item(4, 14).Value = 1.001574323867359
$ws.Cells.Item(4, 15).Value = 1.000590371197722
$ws.Cells.Item(4, 16).Value = 0.999114442882171
$ws.Cells.Item(4, 17).Value = 1.000041472393445
$ws.Cells.Item(4, 18).Value = 0.9999344032105671
$ws.Cells.Item(4, 19).Value = 0.9992404864511704
$ws.Cells.Item(4, 20).Value = 0.9999344032105671
$ws.Cells.Item(4, 21).Value = 0.9998239458052176
$ws.Cells.Item(4, 22).Value = 1.000174021417646
$ws.Cells.Item(4, 23).Value = 0.9998443714206635
$ws.Cells.Item(5, 3).Value = 0.998051177883882
$ws.Cells.Item(5, 4).Value = 0.9902559015513538
$ws.Cells.Item(5, 5).Value = 1.006496075336986
$ws.Cells.Item(5, 6).Value = 0.9902559015513538
$ws.Cells.Item(5, 7).Value = 1.002436014963682
$ws.Cells.Item(5, 8).Value = 1.003616929419082
$ws.Cells.Item(5, 9).Value = 0.9979062274452556
$ws.Cells.Item(5, 10).Value = 1.006496075336986
$ws.Cells.Item(5, 11).Value = 1.006496075336986
$ws.Cells.Item(5, 12).Value = 0.9936643259532137
$ws.Cells.Item(5, 13).Value = 1.002436014963682
$ws.Cells.Item(5, 14).Value = 1.006496075336986
$ws.Cells.Item(5, 15).Value = 1.002436014963682
$ws.Cells.Item(5, 16).Value = 0.996345958257518
$ws.Cells.Item(5, 17).Value = 1.000171121204469
$ws.Cells.Item(5, 18).Value = 0.9997293306173406
$ws.Cells.Item(5, 19).Value = 0.9968660479867638
$ws.Cells.Item(5, 20).Value = 0.9997293306173406
$ws.Cells.Item(5, 21).Value = 0.9992735548243193
$ws.Cells.Item(5, 22).Value = 1.000718058926853
$ws.Cells.Item(5, 23).Value = 0.9993578334396422
$ws.Cells.Item(6, 3).Value = 0.9992403105220923
$ws.Cells.Item(6, 4).Value = 0.996201552908513
$ws.Cells.Item(6, 5).Value = 1.002532298659935
$ws.Cells.Item(6, 6).Value = 0.996201552908513
$ws.Cells.Item(6, 7).Value = 1.000949610662093
$ws.Cells.Item(6, 8).Value = 1.001409950250393
$ws.Cells.Item(6, 9).Value = 0.9991838054932423
$ws.Cells.Item(6, 10).Value = 1.002532298659935
$ws.Cells.Item(6, 11).Value = 1.002532298659935
$ws.Cells.Item(6, 12).Value = 0.9975302287585009
$ws.Cells.Item(6, 13).Value = 1.000949610662093
$ws.Cells.Item(6, 14).Value = 1.002532298659935
$ws.Cells.Item(6, 15).Value = 1.000949610662093
$ws.Cells.Item(6, 16).Value = 0.9985755817853033
$ws.Cells.Item(6, 17).Value = 1.000066708077668
$ws.Cells.Item(6, 18).Value = 0.9998944874101805
$ws.Cells.Item(6, 19).Value = 0.9987783230212829
$ws.Cells.Item(6, 20).Value = 0.9998944874101804
$ws.Cells.Item(6, 21).Value = 0.9997168169309458
$ws.Cells.Item(6, 22).Value = 1.000279913276743
$ws.Cells.Item(6, 23).Value = 0.999749670989608
$ws.Cells.Item(7, 3).Value = 0.9993942567363121
$ws.Cells.Item(7, 4).Value = 0.9969712824855896
$ws.Cells.Item(7, 5).Value = 1.002019145187323
$ws.Cells.Item(7, 6).Value = 0.9969712824855896
$ws.Cells.Item(7, 7).Value = 1.000757180086457
$ws.Cells.Item(7, 8).Value = 1.00112423600144
$ws.Cells.Item(7, 9).Value = 0.9993492013040314
$ws.Cells.Item(7, 10).Value = 1.002019145187323
$ws.Cells.Item(7, 11).Value = 1.002019145187323
$ws.Cells.Item(7, 12).Value = 0.9980307111887606
$ws.Cells.Item(7, 13).Value = 1.000757180086457
$ws.Cells.Item(7, 14).Value = 1.002019145187323
$ws.Cells.Item(7, 15).Value = 1.000757180086457
$ws.Cells.Item(7, 16).Value = 0.9988642312860233
$ws.Cells.Item(7, 17).Value = 1.000053190695244
$ws.Cells.Item(7, 18).Value = 0.9999158692531233
$ws.Cells.Item(7, 19).Value = 0.9990258879586927
$ws.Cells.Item(7, 20).Value = 0.9999158692531233
$ws.Cells.Item(7, 21).Value = 0.9997742022658502
$ws.Cells.Item(7, 22).Value = 1.000223190850145
$ws.Cells.Item(7, 23).Value = 0.9998003991345463
$ws.Cells.Item(8, 3).Value = 0.9980390964658836
$ws.Cells.Item(8, 4).Value = 0.9901954815167299
$ws.Cells.Item(8, 5).Value = 1.006536345882048
$ws.Cells.Item(8, 6).Value = 0.9901954815167299
$ws.Cells.Item(8, 7).Value = 1.002451148303192
$ws.Cells.Item(8, 8).Value = 1.003639368733281
$ws.Cells.Item(8, 9).Value = 0.9978932415120436
$ws.Cells.Item(8, 10).Value = 1.006536345882048
$ws.Cells.Item(8, 11).Value = 1.006536345882048
$ws.Cells.Item(8, 12).Value = 0.9936250460252111
$ws.Cells.Item(8, 13).Value = 1.002451148303192
$ws.Cells.Item(8, 14).Value = 1.006536345882048
$ws.Cells.Item(8, 15).Value = 1.002451148303192
$ws.Cells.Item(8, 16).Value = 0.9963233149099611
$ws.Cells.Item(8, 17).Value = 1.000172194907618
$ws.Cells.Item(8, 18).Value = 0.9997276585673233
$ws.Cells.Item(8, 19).Value = 0.9968466237773219
$ws.Cells.Item(8, 20).Value = 0.9997276585673234
$ws.Cells.Item(8, 21).Value = 0.9992690543035034
$ws.Cells.Item(8, 22).Value = 1.000722512619212
$ws.Cells.Item(8, 23).Value = 0.9993538595926978
$ws.Cells.Item(9, 3).Value = 0.9999328220380707
$ws.Cells.Item(9, 4).Value = 0.9996641066691648
$ws.Cells.Item(9, 5).Value = 1.000223928850732
$ws.Cells.Item(9, 6).Value = 0.9996641066691648
$ws.Cells.Item(9, 7).Value = 1.000083972657827
$ws.Cells.Item(9, 8).Value = 1.000124681402637
$ws.Cells.Item(9, 9).Value = 0.9999278253086085
$ws.Cells.Item(9, 10).Value = 1.000223928850732
$ws.Cells.Item(9, 11).Value = 1.000223928850732
$ws.Cells.Item(9, 12).Value = 0.9997816013186381
$ws.Cells.Item(9, 13).Value = 1.000083972657827
$ws.Cells.Item(9, 14).Value = 1.000223928850732
$ws.Cells.Item(9, 15).Value = 1.000083972657827
$ws.Cells.Item(9, 16).Value = 0.999874039663496
$ws.Cells.Item(9, 17).Value = 1.000005898983218
$ws.Cells.Item(9, 18).Value = 0.9999906693925747
$ws.Cells.Item(9, 19).Value = 0.9998919682118669
$ws.Cells.Item(9, 20).Value = 0.9999906693925746
$ws.Cells.Item(9, 21).Value = 0.9999749583715831
$ws.Cells.Item(9, 22).Value = 1.000024752467413
$ws.Cells.Item(9, 23).Value = 0.9999778638629382
$ws.Cells.Item(10, 3).Value = 0.9999951145350165
$ws.Cells.Item(10, 4).Value = 0.9999755736196053
$ws.Cells.Item(10, 5).Value = 1.000016284378513
$ws.Cells.Item(10, 6).Value = 0.9999755736196053
$ws.Cells.Item(10, 7).Value = 1.000006105717792
$ws.Cells.Item(10, 8).Value = 1.000009067429393
$ws.Cells.Item(10, 9).Value = 0.999994752035404
$ws.Cells.Item(10, 10).Value = 1.000016284378513
$ws.Cells.Item(10, 11).Value = 1.000016284378513
$ws.Cells.Item(10, 12).Value = 0.999984118920814
$ws.Cells.Item(10, 13).Value = 1.000006105717792
$ws.Cells.Item(10, 14).Value = 1.000016284378513
$ws.Cells.Item(10, 15).Value = 1.000006105717792
$ws.Cells.Item(10, 16).Value = 0.9999908396686987
$ws.Cells.Item(10, 17).Value = 1.000000428876598
$ws.Cells.Item(10, 18).Value = 0.9999993212386368
$ws.Cells.Item(10, 19).Value = 0.9999921437909339
$ws.Cells.Item(10, 20).Value = 0.9999993212386368
$ws.Cells.Item(10, 21).Value = 0.9999981789378285
$ws.Cells.Item(10, 22).Value = 1.000001800025965
$ws.Cells.Item(10, 23).Value = 0.9999983902942912
$ws.Cells.Item(11, 3).Value = 0.9998921545033813
$ws.Cells.Item(11, 4).Value = 0.9994607752654463
$ws.Cells.Item(11, 5).Value = 1.000359483073277
$ws.Cells.Item(11, 6).Value = 0.9994607752654463
$ws.Cells.Item(11, 7).Value = 1.00013480587907
$ws.Cells.Item(11, 8).Value = 1.000200156810134
$ws.Cells.Item(11, 9).Value = 0.9998841329181529
$ws.Cells.Item(11, 10).Value = 1.000359483073277
$ws.Cells.Item(11, 11).Value = 1.000359483073277
$ws.Cells.Item(11, 12).Value = 0.9996493947247888
$ws.Cells.Item(11, 13).Value = 1.00013480587907
$ws.Cells.Item(11, 14).Value = 1.000359483073277
$ws.Cells.Item(11, 15).Value = 1.00013480587907
$ws.Cells.Item(11, 16).Value = 0.9997977905722584
$ws.Cells.Item(11, 17).Value = 1.000009469398612
$ws.Cells.Item(11, 18).Value = 0.9999850214059314
$ws.Cells.Item(11, 19).Value = 0.9998265713542231
$ws.Cells.Item(11, 20).Value = 0.9999850214059314
$ws.Cells.Item(11, 21).Value = 0.9999597992839868
$ws.Cells.Item(11, 22).Value = 1.000039736041845
$ws.Cells.Item(11, 23).Value = 0.9999644636316652
$ws.Cells.Item(12, 3).Value = 0.9980099659869232
$ws.Cells.Item(12, 4).Value = 0.9900498284764824
$ws.Cells.Item(12, 5).Value = 1.006633448352998
$ws.Cells.Item(12, 6).Value = 0.9900498284764824
$ws.Cells.Item(12, 7).Value = 1.002487562034373
$ws.Cells.Item(12, 8).Value = 1.003693434601515
$ws.Cells.Item(12, 9).Value = 0.9978619438379697
$ws.Cells.Item(12, 10).Value = 1.006633448352998
$ws.Cells.Item(12, 11).Value = 1.006633448352998
$ws.Cells.Item(12, 12).Value = 0.9935303408479749
$ws.Cells.Item(12, 13).Value = 1.002487562034373
$ws.Cells.Item(12, 14).Value = 1.006633448352998
$ws.Cells.Item(12, 15).Value = 1.002487562034373
$ws.Cells.Item(12, 16).Value = 0.9962686952554276
$ws.Cells.Item(12, 17).Value = 1.000174752936171
$ws.Cells.Item(12, 18).Value = 0.9997236129546175
$ws.Cells.Item(12, 19).Value = 0.9967997781162751
$ws.Cells.Item(12, 20).Value = 0.9997236129546175
$ws.Cells.Item(12, 21).Value = 0.9992581956754556
$ws.Cells.Item(12, 22).Value = 1.000733246210964
$ws.Cells.Item(12, 23).Value = 0.999344260771576
$ws.Cells.Item(13, 3).Value = 0.9996136300225946
$ws.Cells.Item(13, 4).Value = 0.9980681522933327
$ws.Cells.Item(13, 5).Value = 1.001287898614653
$ws.Cells.Item(13, 6).Value = 0.9980681522933327
$ws.Cells.Item(13, 7).Value = 1.000482959917478
$ws.Cells.Item(13, 8).Value = 1.000717082759434
$ws.Cells.Item(13, 9).Value = 0.9995848929847293
$ws.Cells.Item(13, 10).Value = 1.001287898614653
$ws.Cells.Item(13, 11).Value = 1.001287898614653
$ws.Cells.Item(13, 12).Value = 0.9987439030728202
$ws.Cells.Item(13, 13).Value = 1.000482959917478
$ws.Cells.Item(13, 14).Value = 1.001287898614653
$ws.Cells.Item(13, 15).Value = 1.000482959917478
$ws.Cells.Item(13, 16).Value = 0.9992755561054052
$ws.Cells.Item(13, 17).Value = 1.000033926451104
$ws.Cells.Item(13, 18).Value = 0.999946336941821
$ws.Cells.Item(13, 19).Value = 0.9993786683985132
$ws.Cells.Item(13, 20).Value = 0.999946336941821
$ws.Cells.Item(13, 21).Value = 0.999855975952548
$ws.Cells.Item(13, 22).Value = 1.000142360484969
$ws.Cells.Item(13, 23).Value = 0.9998726849478149
$ws.Cells.Item(14, 3).Value = 0.9993950757999988
$ws.Cells.Item(14, 4).Value = 0.9969753752631567
$ws.Cells.Item(14, 5).Value = 1.002016417684213
$ws.Cells.Item(14, 6).Value = 0.9969753752631567
$ws.Cells.Item(14, 7).Value = 1.000756158147369
$ws.Cells.Item(14, 8).Value = 1.001122719021054
$ws.Cells.Item(14, 9).Value = 0.9993500805789479
$ws.Cells.Item(14, 10).Value = 1.002016417684213
$ws.Cells.Item(14, 11).Value = 1.002016417684213
$ws.Cells.Item(14, 12).Value = 0.9980333708210524
$ws.Cells.Item(14, 13).Value = 1.000756158147369
$ws.Cells.Item(14, 14).Value = 1.002016417684213
$ws.Cells.Item(14, 15).Value = 1.000756158147369
$ws.Cells.Item(14, 16).Value = 0.9988657667052631
$ws.Cells.Item(14, 17).Value = 1.000053119363159
$ws.Cells.Item(14, 18).Value = 0.9999159836982464
$ws.Cells.Item(14, 19).Value = 0.999027204663158
$ws.Cells.Item(14, 20).Value = 0.9999159836982464
$ws.Cells.Item(14, 21).Value = 0.9997745079184217
$ws.Cells.Item(14, 22).Value = 1.00022288987158
$ws.Cells.Item(14, 23).Value = 0.9998006694328953
$ws.Cells.Item(15, 3).Value = 1.000391390152756
$ws.Cells.Item(15, 4).Value = 1.001956946792648
$ws.Cells.Item(15, 5).Value = 0.9986953714659608
$ws.Cells.Item(15, 6).Value = 1.001956946792648
$ws.Cells.Item(15, 7).Value = 0.9995107639784528
$ws.Cells.Item(15, 8).Value = 0.9992736015494776
$ws.Cells.Item(15, 9).Value = 1.000420501303533
$ws.Cells.Item(15, 10).Value = 0.9986953714659608
$ws.Cells.Item(15, 11).Value = 0.9986953714659608
$ws.Cells.Item(15, 12).Value = 1.001272417166038
$ws.Cells.Item(15, 13).Value = 0.9995107639784528
$ws.Cells.Item(15, 14).Value = 0.9986953714659608
$ws.Cells.Item(15, 15).Value = 0.9995107639784528
$ws.Cells.Item(15, 16).Value = 1.00073385538555
$ws.Cells.Item(15, 17).Value = 0.9999656326409928
$ws.Cells.Item(15, 18).Value = 1.000054360745687
$ws.Cells.Item(15, 19).Value = 1.000629404024878
$ws.Cells.Item(15, 20).Value = 1.000054360745687
$ws.Cells.Item(15, 21).Value = 1.000145895885149
$ws.Cells.Item(15, 22).Value = 0.9998557910013111
$ws.Cells.Item(15, 23).Value = 1.000128969548415
$ws.Cells.Item(16, 3).Value = 0.9965617900000008
$ws.Cells.Item(16, 4).Value = 0.98280895
$ws.Cells.Item(16, 5).Value = 1.0114607
$ws.Cells.Item(16, 6).Value = 0.98280895
$ws.Cells.Item(16, 7).Value = 1.004297800000002
$ws.Cells.Item(16, 8).Value = 1.006381200000001
$ws.Cells.Item(16, 9).Value = 0.9963060500000007
$ws.Cells.Item(16, 10).Value = 1.0114607
$ws.Cells.Item(16, 11).Value = 1.0114607
$ws.Cells.Item(16, 12).Value = 0.9888222800000002
$ws.Cells.Item(16, 13).Value = 1.004297800000002
$ws.Cells.Item(16, 14).Value = 1.0114607
$ws.Cells.Item(16, 15).Value = 1.004297800000002
$ws.Cells.Item(16, 16).Value = 0.993553375000001
$ws.Cells.Item(16, 17).Value = 1.000301925000001
$ws.Cells.Item(16, 18).Value = 0.999522483333334
$ws.Cells.Item(16, 19).Value = 0.9944709333333343
$ws.Cells.Item(16, 20).Value = 0.999522483333334
$ws.Cells.Item(16, 21).Value = 0.9987183750000007
$ws.Cells.Item(16, 22).Value = 1.001266840000001
$ws.Cells.Item(16, 23).Value = 0.9988670712500008
$ws.Cells.Item(17, 3).Value = 0.9944594800000001
$ws.Cells.Item(17, 4).Value = 0.9722974
$ws.Cells.Item(17, 5).Value = 1.0184684
$ws.Cells.Item(17, 6).Value = 0.9722974
$ws.Cells.Item(17, 7).Value = 1.0069256
$ws.Cells.Item(17, 8).Value = 1.010283
$ws.Cells.Item(17, 9).Value = 0.99404738
$ws.Cells.Item(17, 10).Value = 1.0184684
$ws.Cells.Item(17, 11).Value = 1.0184684
$ws.Cells.Item(17, 12).Value = 0.98198761
$ws.Cells.Item(17, 13).Value = 1.0069256
$ws.Cells.Item(17, 14).Value = 1.0184684
$ws.Cells.Item(17, 15).Value = 1.0069256
$ws.Cells.Item(17, 16).Value = 0.9896115
$ws.Cells.Item(17, 17).Value = 1.00048649
$ws.Cells.Item(17, 18).Value = 0.9992304666666666
$ws.Cells.Item(17, 19).Value = 0.9910901266666666
$ws.Cells.Item(17, 20).Value = 0.9992304666666666
$ws.Cells.Item(17, 21).Value = 0.997934695
$ws.Cells.Item(17, 22).Value = 1.002041436
$ws.Cells.Item(17, 23).Value = 0.99817430875
$ws.Cells.Item(18, 3).Value = 0.99466084
$ws.Cells.Item(18, 4).Value = 0.9733042199999999
$ws.Cells.Item(18, 5).Value = 1.0177972
$ws.Cells.Item(18, 6).Value = 0.9733042199999999
$ws.Cells.Item(18, 7).Value = 1.0066739
$ws.Cells.Item(18, 8).Value = 1.0099092
$ws.Cells.Item(18, 9).Value = 0.99426372
$ws.Cells.Item(18, 10).Value = 1.0177972
$ws.Cells.Item(18, 11).Value = 1.0177972
$ws.Cells.Item(18, 12).Value = 0.9826422500000001
$ws.Cells.Item(18, 13).Value = 1.0066739
$ws.Cells.Item(18, 14).Value = 1.0177972
$ws.Cells.Item(18, 15).Value = 1.0066739
$ws.Cells.Item(18, 16).Value = 0.98998906
$ws.Cells.Item(18, 17).Value = 1.00046881
$ws.Cells.Item(18, 18).Value = 0.9992584399999999
$ws.Cells.Item(18, 19).Value = 0.9914139466666666
$ws.Cells.Item(18, 20).Value = 0.9992584399999999
$ws.Cells.Item(18, 21).Value = 0.99800976
$ws.Cells.Item(18, 22).Value = 1.001967248
$ws.Cells.Item(18, 23).Value = 0.9982406537499999
$ws.Cells.Item(19, 3).Value = 0.99860256
$ws.Cells.Item(19, 4).Value = 0.99301282
$ws.Cells.Item(19, 5).Value = 1.0046581
$ws.Cells.Item(19, 6).Value = 0.99301282
$ws.Cells.Item(19, 7).Value = 1.0017468
$ws.Cells.Item(19, 8).Value = 1.0025936
$ws.Cells.Item(19, 9).Value = 0.99849862
$ws.Cells.Item(19, 10).Value = 1.0046581
$ws.Cells.Item(19, 11).Value = 1.0046581
$ws.Cells.Item(19, 12).Value = 0.9954568900000001
$ws.Cells.Item(19, 13).Value = 1.0017468
$ws.Cells.Item(19, 14).Value = 1.0046581
$ws.Cells.Item(19, 15).Value = 1.0017468
$ws.Cells.Item(19, 16).Value = 0.99737981
$ws.Cells.Item(19, 17).Value = 1.00012271
$ws.Cells.Item(19, 18).Value = 0.9998059066666668
$ws.Cells.Item(19, 19).Value = 0.9977527466666666
$ws.Cells.Item(19, 20).Value = 0.9998059066666668
$ws.Cells.Item(19, 21).Value = 0.999479085
$ws.Cells.Item(19, 22).Value = 1.000514888
$ws.Cells.Item(19, 23).Value = 0.99953952375
$ws.Cells.Item(20, 3).Value = 0.9986524432876712
$ws.Cells.Item(20, 4).Value = 0.9932621802739727
$ws.Cells.Item(20, 5).Value = 1.004491875616438
$ws.Cells.Item(20, 6).Value = 0.9932621802739727
$ws.Cells.Item(20, 7).Value = 1.001684454794521
$ws.Cells.Item(20, 8).Value = 1.002501016986301
$ws.Cells.Item(20, 9).Value = 0.9985522057534247
$ws.Cells.Item(20, 10).Value = 1.004491875616438
$ws.Cells.Item(20, 11).Value = 1.004491875616438
$ws.Cells.Item(20, 12).Value = 0.995619035205479
$ws.Cells.Item(20, 13).Value = 1.001684454794521
$ws.Cells.Item(20, 14).Value = 1.004491875616438
$ws.Cells.Item(20, 15).Value = 1.001684454794521
$ws.Cells.Item(20, 16).Value = 0.9974733175342467
$ws.Cells.Item(20, 17).Value = 1.000118330273973
$ws.Cells.Item(20, 18).Value = 0.9998128368949772
$ws.Cells.Item(20, 19).Value = 0.9978329469406394
$ws.Cells.Item(20, 20).Value = 0.9998128368949772
$ws.Cells.Item(20, 21).Value = 0.9994976791095891
$ws.Cells.Item(20, 22).Value = 1.000496518410959
$ws.Cells.Item(20, 23).Value = 0.999555958339041
$ws.Cells.Item(21, 3).Value = 0.997800675263158
$ws.Cells.Item(21, 4).Value = 0.9890033726315791
$ws.Cells.Item(21, 5).Value = 1.007331088947368
$ws.Cells.Item(21, 6).Value = 0.9890033726315791
$ws.Cells.Item(21, 7).Value = 1.002749164736842
$ws.Cells.Item(21, 8).Value = 1.00408185368421
$ws.Cells.Item(21, 9).Value = 0.9976370826315789
$ws.Cells.Item(21, 10).Value = 1.007331088947368
$ws.Cells.Item(21, 11).Value = 1.007331088947368
$ws.Cells.Item(21, 12).Value = 0.9928499273684211
$ws.Cells.Item(21, 13).Value = 1.002749164736842
$ws.Cells.Item(21, 14).Value = 1.007331088947368
$ws.Cells.Item(21, 15).Value = 1.002749164736842
$ws.Cells.Item(21, 16).Value = 0.9958762686842106
$ws.Cells.Item(21, 17).Value = 1.000193123684211
$ws.Cells.Item(21, 18).Value = 0.9996945421052632
$ws.Cells.Item(21, 19).Value = 0.9964632066666668
$ws.Cells.Item(21, 20).Value = 0.9996945421052633
$ws.Cells.Item(21, 21).Value = 0.9991801772368423
$ws.Cells.Item(21, 22).Value = 1.000810359578947
$ws.Cells.Item(21, 23).Value = 0.99927529125
$ws.Cells.Item(22, 3).Value = 0.9978669105263159
$ws.Cells.Item(22, 4).Value = 0.9893345631578947
$ws.Cells.Item(22, 5).Value = 1.007110298947368
$ws.Cells.Item(22, 6).Value = 0.9893345631578947
$ws.Cells.Item(22, 7).Value = 1.002666350526316
$ws.Cells.Item(22, 8).Value = 1.003958923684211
$ws.Cells.Item(22, 9).Value = 0.9977082536842106
$ws.Cells.Item(22, 10).Value = 1.007110298947368
$ws.Cells.Item(22, 11).Value = 1.007110298947368
$ws.Cells.Item(22, 12).Value = 0.9930652668421055
$ws.Cells.Item(22, 13).Value = 1.002666350526316
$ws.Cells.Item(22, 14).Value = 1.007110298947368
$ws.Cells.Item(22, 15).Value = 1.002666350526316
$ws.Cells.Item(22, 16).Value = 0.9960004568421053
$ws.Cells.Item(22, 17).Value = 1.000187302105263
$ws.Cells.Item(22, 18).Value = 0.9997037375438597
$ws.Cells.Item(22, 19).Value = 0.9965697224561403
$ws.Cells.Item(22, 20).Value = 0.9997037375438597
$ws.Cells.Item(22, 21).Value = 0.9992048665789475
$ws.Cells.Item(22, 22).Value = 1.000785953052632
$ws.Cells.Item(22, 23).Value = 0.9992971147368421
$ws.Cells.Item(23, 3).Value = 1.00000325123103
$ws.Cells.Item(23, 4).Value = 1.000016246698306
$ws.Cells.Item(23, 5).Value = 0.9999891729394074
$ws.Cells.Item(23, 6).Value = 1.000016246698306
$ws.Cells.Item(23, 7).Value = 0.9999959388576023
$ws.Cells.Item(23, 8).Value = 0.9999939776097011
$ws.Cells.Item(23, 9).Value = 1.000003487505612
$ws.Cells.Item(23, 10).Value = 0.9999891729394074
$ws.Cells.Item(23, 11).Value = 0.9999891729394074
$ws.Cells.Item(23, 12).Value = 1.000010562391112
$ws.Cells.Item(23, 13).Value = 0.9999959388576023
$ws.Cells.Item(23, 14).Value = 0.9999891729394074
$ws.Cells.Item(23, 15).Value = 0.9999959388576023
$ws.Cells.Item(23, 16).Value = 1.000006092777954
$ws.Cells.Item(23, 17).Value = 0.999999713181607
$ws.Cells.Item(23, 18).Value = 1.000000452831772
$ws.Cells.Item(23, 19).Value = 1.00000522435384
$ws.Cells.Item(23, 20).Value = 1.000000452831772
$ws.Cells.Item(23, 21).Value = 1.000001211500232
$ws.Cells.Item(23, 22).Value = 0.9999988037880667
$ws.Cells.Item(23, 23).Value = 1.000001072011297
$ws.Cells.Item(24, 3).Value = 0.9999991544416285
$ws.Cells.Item(24, 4).Value = 0.9999957527114924
$ws.Cells.Item(24, 5).Value = 1.000002831313868
$ws.Cells.Item(24, 6).Value = 0.9999957527114924
$ws.Cells.Item(24, 7).Value = 1.000001062561848
$ws.Cells.Item(24, 8).Value = 1.000001579236763
$ws.Cells.Item(24, 9).Value = 0.9999990867814664
$ws.Cells.Item(24, 10).Value = 1.000002831313868
$ws.Cells.Item(24, 11).Value = 1.000002831313868
$ws.Cells.Item(24, 12).Value = 0.9999972423346497
$ws.Cells.Item(24, 13).Value = 1.000001062561848
$ws.Cells.Item(24, 14).Value = 1.000002831313868
$ws.Cells.Item(24, 15).Value = 1.000001062561848
$ws.Cells.Item(24, 16).Value = 0.9999984076366704
$ws.Cells.Item(24, 17).Value = 1.000000074671657
$ws.Cells.Item(24, 18).Value = 0.9999998821957364
$ws.Cells.Item(24, 19).Value = 0.9999986340182691
$ws.Cells.Item(24, 20).Value = 0.9999998821957364
$ws.Cells.Item(24, 21).Value = 0.9999996833421689
$ws.Cells.Item(24, 22).Value = 1.000000312936509
$ws.Cells.Item(24, 23).Value = 0.9999997214929457
$ws.Cells.Item(25, 3).Value = 1.000802506262872
$ws.Cells.Item(25, 4).Value = 1.00401254475068
$ws.Cells.Item(25, 5).Value = 0.9973249710934048
$ws.Cells.Item(25, 6).Value = 1.00401254475068
$ws.Cells.Item(25, 7).Value = 0.9989968660173206
$ws.Cells.Item(25, 8).Value = 0.9985105793061635
$ws.Cells.Item(25, 9).Value = 1.000862207233465
$ws.Cells.Item(25, 10).Value = 0.9973249710934048
$ws.Cells.Item(25, 11).Value = 0.9973249710934048
$ws.Cells.Item(25, 12).Value = 1.002608983010344
$ws.Cells.Item(25, 13).Value = 0.9989968660173206
$ws.Cells.Item(25, 14).Value = 0.9973249710934048
$ws.Cells.Item(25, 15).Value = 0.9989968660173206
$ws.Cells.Item(25, 16).Value = 1.001504705384
$ws.Cells.Item(25, 17).Value = 0.9999295366253929
$ws.Cells.Item(25, 18).Value = 1.000111460620468
$ws.Cells.Item(25, 19).Value = 1.001290539333822
$ws.Cells.Item(25, 20).Value = 1.000111460620468
$ws.Cells.Item(25, 21).Value = 1.000299147273718
$ws.Cells.Item(25, 22).Value = 0.9997043120376551
$ws.Cells.Item(25, 23).Value = 1.000264440461446
$ws.Cells.Item(26, 3).Value = 1.000208009063202
$ws.Cells.Item(26, 4).Value = 1.001040039616234
$ws.Cells.Item(26, 5).Value = 0.9993066421812604
$ws.Cells.Item(26, 6).Value = 1.001040039616234
$ws.Cells.Item(26, 7).Value = 0.9997399888807381
$ws.Cells.Item(26, 8).Value = 0.9996139498440455
$ws.Cells.Item(26, 9).Value = 1.00022348233041
$ws.Cells.Item(26, 10).Value = 0.9993066421812604
$ws.Cells.Item(26, 11).Value = 0.9993066421812604
$ws.Cells.Item(26, 12).Value = 1.000676236418479
$ws.Cells.Item(26, 13).Value = 0.9997399888807381
$ws.Cells.Item(26, 14).Value = 0.9993066421812604
$ws.Cells.Item(26, 15).Value = 0.9997399888807381
$ws.Cells.Item(26, 16).Value = 1.000390014248486
$ws.Cells.Item(26, 17).Value = 0.9999817356055742
$ws.Cells.Item(26, 18).Value = 1.000028890226077
$ws.Cells.Item(26, 19).Value = 1.000334503609127
$ws.Cells.Item(26, 20).Value = 1.000028890226077
$ws.Cells.Item(26, 21).Value = 1.000077538252161
$ws.Cells.Item(26, 22).Value = 0.9999233590379806
$ws.Cells.Item(26, 23).Value = 1.000068542151888
$ws.Cells.Item(27, 3).Value = 1.000194020052813
$ws.Cells.Item(27, 4).Value = 1.000970044224112
$ws.Cells.Item(27, 5).Value = 0.9993533042348421
$ws.Cells.Item(27, 6).Value = 1.000970044224112
$ws.Cells.Item(27, 7).Value = 0.9997574869348043
$ws.Cells.Item(27, 8).Value = 0.9996399257986494
$ws.Cells.Item(27, 9).Value = 1.000208433244725
$ws.Cells.Item(27, 10).Value = 0.9993533042348421
$ws.Cells.Item(27, 11).Value = 0.9993533042348421
$ws.Cells.Item(27, 12).Value = 1.000630731069611
$ws.Cells.Item(27, 13).Value = 0.9997574869348043
$ws.Cells.Item(27, 14).Value = 0.9993533042348421
$ws.Cells.Item(27, 15).Value = 0.9997574869348043
$ws.Cells.Item(27, 16).Value = 1.000363765579458
$ws.Cells.Item(27, 17).Value = 0.9999829600897648
$ws.Cells.Item(27, 18).Value = 1.000026945131253
$ws.Cells.Item(27, 19).Value = 1.000311988134547
$ws.Cells.Item(27, 20).Value = 1.000026945131253
$ws.Cells.Item(27, 21).Value = 1.000072317159621
$ws.Cells.Item(27, 22).Value = 0.999928514574665
$ws.Cells.Item(27, 23).Value = 1.000063929061795
$ws.Cells.Item(28, 3).Value = 1.000050163014436
$ws.Cells.Item(28, 4).Value = 1.000250803452359
$ws.Cells.Item(28, 5).Value = 0.9998328016334009
$ws.Cells.Item(28, 6).Value = 1.000250803452359
$ws.Cells.Item(28, 7).Value = 0.9999372997620529
$ws.Cells.Item(28, 8).Value = 0.9999069072746795
$ws.Cells.Item(28, 9).Value = 1.000053888395534
$ws.Cells.Item(28, 10).Value = 0.9998328016334009
$ws.Cells.Item(28, 11).Value = 0.9998328016334009
$ws.Cells.Item(28, 12).Value = 1.000163072069299
$ws.Cells.Item(28, 13).Value = 0.9999372997620529
$ws.Cells.Item(28, 14).Value = 0.9998328016334009
$ws.Cells.Item(28, 15).Value = 0.9999372997620529
$ws.Cells.Item(28, 16).Value = 1.000094051607206
$ws.Cells.Item(28, 17).Value = 0.9999955940787937
$ws.Cells.Item(28, 18).Value = 1.000006968282604
$ws.Cells.Item(28, 19).Value = 1.000080663869982
$ws.Cells.Item(28, 20).Value = 1.000006968282604
$ws.Cells.Item(28, 21).Value = 1.000018698310837
$ws.Cells.Item(28, 22).Value = 0.9999815189753495
$ws.Cells.Item(28, 23).Value = 1.000016529420477
$ws.Cells.Item(29, 3).Value = 0.9997950679831376
$ws.Cells.Item(29, 4).Value = 0.9989754122851008
$ws.Cells.Item(29, 5).Value = 1.000683055130593
$ws.Cells.Item(29, 6).Value = 0.9989754122851008
$ws.Cells.Item(29, 7).Value = 1.000256150288184
$ws.Cells.Item(29, 8).Value = 1.000380313718002
$ws.Cells.Item(29, 9).Value = 0.999779831212017
$ws.Cells.Item(29, 10).Value = 1.000683055130593
$ws.Cells.Item(29, 11).Value = 1.000683055130593
$ws.Cells.Item(29, 12).Value = 0.999333820255321
$ws.Cells.Item(29, 13).Value = 1.000256150288184
$ws.Cells.Item(29, 14).Value = 1.000683055130593
$ws.Cells.Item(29, 15).Value = 1.000256150288184
$ws.Cells.Item(29, 16).Value = 0.9996157812866424
$ws.Cells.Item(29, 17).Value = 1.000017990750101
$ws.Cells.Item(29, 18).Value = 0.9999715392346259
$ws.Cells.Item(29, 19).Value = 0.9996704645951006
$ws.Cells.Item(29, 20).Value = 0.9999715392346259
$ws.Cells.Item(29, 21).Value = 0.9999236122289736
$ws.Cells.Item(29, 22).Value = 1.000075500809297
$ws.Cells.Item(29, 23).Value = 0.9999324751450674
$ws.Cells.Item(30, 3).Value = 0.9997785210774021
$ws.Cells.Item(30, 4).Value = 0.9988925978481382
$ws.Cells.Item(30, 5).Value = 1.000738268854098
$ws.Cells.Item(30, 6).Value = 0.9988925978481382
$ws.Cells.Item(30, 7).Value = 1.000276855181368
$ws.Cells.Item(30, 8).Value = 1.000411059032649
$ws.Cells.Item(30, 9).Value = 0.999762040931674
$ws.Cells.Item(30, 10).Value = 1.000738268854098
$ws.Cells.Item(30, 11).Value = 1.000738268854098
$ws.Cells.Item(30, 12).Value = 0.9992799566510927
$ws.Cells.Item(30, 13).Value = 1.000276855181368
$ws.Cells.Item(30, 14).Value = 1.000738268854098
$ws.Cells.Item(30, 15).Value = 1.000276855181368
$ws.Cells.Item(30, 16).Value = 0.9995847265147532
$ws.Cells.Item(30, 17).Value = 1.000019448056521
$ws.Cells.Item(30, 18).Value = 0.9999692406278681
$ws.Cells.Item(30, 19).Value = 0.9996438313203934
$ws.Cells.Item(30, 20).Value = 0.9999692406278681
$ws.Cells.Item(30, 21).Value = 0.9999174407038196
$ws.Cells.Item(30, 22).Value = 1.000081606333875
$ws.Cells.Item(30, 23).Value = 0.9999270193447238
$ws.Cells.Item(31, 3).Value = 1.000407080863157
$ws.Cells.Item(31, 4).Value = 1.002035377164536
$ws.Cells.Item(31, 5).Value = 0.9986430942741782
$ws.Cells.Item(31, 6).Value = 1.002035377164536
$ws.Cells.Item(31, 7).Value = 0.9994911784473592
$ws.Cells.Item(31, 8).Value = 0.9992444917654071
$ws.Cells.Item(31, 9).Value = 1.000437344205353
$ws.Cells.Item(31, 10).Value = 0.9986430942741782
$ws.Cells.Item(31, 11).Value = 0.9986430942741782
$ws.Cells.Item(31, 12).Value = 1.001323404090686
$ws.Cells.Item(31, 13).Value = 0.9994911784473592
$ws.Cells.Item(31, 14).Value = 0.9986430942741782
$ws.Cells.Item(31, 15).Value = 0.9994911784473592
$ws.Cells.Item(31, 16).Value = 1.000763277805948
$ws.Cells.Item(31, 17).Value = 0.999964261326356
$ws.Cells.Item(31, 18).Value = 1.000056549962024
$ws.Cells.Item(31, 19).Value = 1.000654633272416
$ws.Cells.Item(31, 20).Value = 1.000056549962024
$ws.Cells.Item(31, 21).Value = 1.000151748522857
$ws.Cells.Item(31, 22).Value = 0.9998500176731209
$ws.Cells.Item(31, 23).Value = 1.000134143657254
